$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 102, shifting rows 102:226 down to 103:227
$ws.Rows("102").Insert()

# Populate the newly inserted row 102 with the new record
$ws.Range("A102").Value = 3
$ws.Range("B102").Value = "Femacal de La Calera"
$ws.Range("C102").Value = "Coquimbo"
$ws.Range("D102").Value = 44546
$ws.Range("E102").Value = 5
$ws.Range("F102").Value = 100112039
$ws.Range("G102").Value = "Ciboulette"
$ws.Range("H102").Value = "Sin especificar"
$ws.Range("I102").Value = "Primera"
$ws.Range("J102").Value = 180
$ws.Range("K102").Value = 1500
$ws.Range("L102").Value = 1500
$ws.Range("M102").Value = 1500
$ws.Range("N102").Value = "$/docena de atados"
$ws.Range("O102").Value = "Provincia de Quillota"
$ws.Range("P102").Value = 500
$ws.Range("Q102").Value = 3
$ws.Range("R102").Value = "Hortaliza"
